# httk-benchmarks.xlsx — add the "2.5.0" benchmark row (models 3comp2 and
# sumclearances) to Table1 on Sheet1, mirroring what the F26/J26-derived
# chart series caches already show for the new 26th category point.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend Table1 by one row — this grows ref/autoFilter from A1:R26 to
# A1:R27 and keeps the worksheet's used-range dimension in sync.
$tbl = $ws.ListObjects.Item("Table1")
$newRow = $tbl.ListRows.Add()
$newRange = $newRow.Range

# Match the left-aligned style ("s=1") used by every other data row.
$newRange.HorizontalAlignment = -4131   # xlLeft

# Row 27 values (same layout as the rest of the table):
# Version, N.steadystate, calc_analytic.units, calc_mc.units,
# solve_pbtk.units, RMSLE.Wetmore, N.Wetmore, RMSLE.noMC, N.noMC,
# RMSLE.InVivoCss, N.InVivoCss, RMSLE.InVivoAUC, N.InVivoAUC,
# RMSLE.InVivoCmax, N.InVivoCmax, RMSLE.TissuePC, N.TissuePC, Notes
$ws.Range("A27").Value2 = "2.5.0"
$ws.Range("B27").Value2 = 1021
$ws.Range("C27").Value2 = 1
$ws.Range("D27").Value2 = 1
$ws.Range("E27").Value2 = 0.99990000000000001
$ws.Range("F27").Value2 = 0.94769999999999999
$ws.Range("G27").Value2 = 353
$ws.Range("H27").Value2 = 0.27160000000000001
$ws.Range("I27").Value2 = 353
$ws.Range("J27").Value2 = 1.508
$ws.Range("K27").Value2 = 36
$ws.Range("L27").Value2 = 0.9698
$ws.Range("M27").Value2 = 80
$ws.Range("N27").Value2 = 1.1319999999999999
$ws.Range("O27").Value2 = 80
$ws.Range("P27").Value2 = 0.64659999999999995
$ws.Range("Q27").Value2 = 863
$ws.Range("R27").Value2 = "Added models 3comp2 and sumclearances"

# Leave the same kind of selection Excel would after typing into the new
# row's last (Notes) cell.
[void]$ws.Activate()
[void]$ws.Range("R27").Select()
